# added 4wk low sales check
# Updates forecast figures on the "Forecast Comparison" sheet (MyForecast,
# Inventory Coverage, Seasonality Index) and the corresponding roll-up
# figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$fc = $wb.Worksheets.Item("Forecast Comparison")
$sm = $wb.Worksheets.Item("Summary")

# Row => (MyForecast, Inventory Coverage, Seasonality Index)
$fc.Range("D2").Value  = 154
$fc.Range("H2").Value  = 12.92
$fc.Range("L2").Value  = 1.02

$fc.Range("D3").Value  = 155
$fc.Range("H3").Value  = 11.88
$fc.Range("L3").Value  = 0.88

$fc.Range("D4").Value  = 155
$fc.Range("H4").Value  = 10.85
$fc.Range("L4").Value  = 0.91

$fc.Range("D5").Value  = 156
$fc.Range("H5").Value  = 9.83
$fc.Range("L5").Value  = 0.92

$fc.Range("D6").Value  = 156
$fc.Range("H6").Value  = 8.800000000000001
$fc.Range("L6").Value  = 1.04

$fc.Range("D7").Value  = 157
$fc.Range("H7").Value  = 7.78
$fc.Range("L7").Value  = 0.8100000000000001

$fc.Range("D8").Value  = 157
$fc.Range("H8").Value  = 6.77
$fc.Range("L8").Value  = 0.99

$fc.Range("D9").Value  = 157
$fc.Range("H9").Value  = 5.75
$fc.Range("L9").Value  = 1.16

$fc.Range("D10").Value = 158
$fc.Range("H10").Value = 4.74
$fc.Range("L10").Value = 0.82

$fc.Range("D11").Value = 158
$fc.Range("H11").Value = 3.73
$fc.Range("L11").Value = 0.9

$fc.Range("D12").Value = 159
$fc.Range("H12").Value = 2.72
$fc.Range("L12").Value = 0.96

$fc.Range("D13").Value = 159
$fc.Range("H13").Value = 1.72
$fc.Range("L13").Value = 0.9399999999999999

$fc.Range("D14").Value = 159
$fc.Range("H14").Value = 0.72
$fc.Range("L14").Value = 0.9399999999999999

$fc.Range("D15").Value = 160
$fc.Range("L15").Value = 1

$fc.Range("D16").Value = 160
$fc.Range("L16").Value = 1.07

$fc.Range("D17").Value = 161
$fc.Range("L17").Value = 0.89

# Summary roll-up figures (kept as text, matching the source file's
# inlineStr representation for this column)
$sm.Range("B9").NumberFormat  = "@"
$sm.Range("B9").Value  = "2529"

$sm.Range("B10").NumberFormat = "@"
$sm.Range("B10").Value = "1252"

$sm.Range("B11").NumberFormat = "@"
$sm.Range("B11").Value = "623"

$sm.Range("B12").NumberFormat = "@"
$sm.Range("B12").Value = "161"

$sm.Range("B14").NumberFormat = "@"
$sm.Range("B14").Value = "155"
